$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.293.03"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.39%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.505.93"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.18%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "570.99"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.58%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "165.09"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.72%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.511"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.24%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.504.15"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.24%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.158"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.17%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.167"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.53%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.356"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.15%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.88"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.32%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.968.69"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "69.188.41"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.32%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000174"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.09%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "24.73"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.88%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.514.11"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.61%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.25"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.99%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.53"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.41%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "348.00"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.16%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.89"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.18%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.98"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.86%  "
$ws.Range("E24").Value = "  -0.07%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "70.14"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.66%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.89"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.10%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.85"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.73%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.654.04"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.46%  "
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0880"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.13%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.80"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.19%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "459.11"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.64%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.23"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -6.04%  "
$ws.Range("E34").Value = "  -2.19%  "
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "157.17"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.26%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.115"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.21%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.07"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.92%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.38"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.88%  "
$ws.Range("E40").Value = "  -0.02%  "
$ws.Range("B41").Value = "PolygonEcosystemToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.316"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.06%  "
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.68"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.87%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.59"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.68%  "
$ws.Range("E44").Value = "  -0.48%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.11"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.76%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.21"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -6.55%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "141.33"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.85%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.46"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.72%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.517"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.73%  "
$ws.Range("E50").Value = "  +0.17%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.576"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.48%  "
